$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 96, shifting existing rows 96:189 down to 97:190
$ws.Rows("96:96").Insert()

# Populate the newly inserted row 96 with the new record
$ws.Range("A96").Value = 9
$ws.Range("B96").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C96").Value = "Metropolitana"
$ws.Range("D96").Value = 44539
$ws.Range("E96").Value = 13
$ws.Range("F96").Value = 100112026
$ws.Range("G96").Value = "Haba"
$ws.Range("H96").Value = "Sin especificar"
$ws.Range("I96").Value = "Primera"
$ws.Range("J96").Value = 61
$ws.Range("K96").Value = 8000
$ws.Range("L96").Value = 9000
$ws.Range("M96").Value = 8492
$ws.Range("N96").Value = "`$/saco 25 kilos"
$ws.Range("O96").Value = "Región del Maule"
$ws.Range("P96").Value = 340
$ws.Range("Q96").Value = 25
$ws.Range("R96").Value = "Hortaliza"
